{"js": "// Apply the two red-colored annotation edits requested by Oliver.\n// Both edits split the leading review run so the leading space stays\n// in the original (uncolored) run while the remaining review text\n// (and, for the second comment, the rest of the paragraph's runs)\n// becomes font color C9211E, keeping any existing highlight.\n\nconst body = context.document.body;\n\n// --- Edit 1: \"... on page 56 we find ...\" comment -----------------\nconst target1 =\n  \"on page 56 we find \\u201coptimization\\u201d and \\u201cneighboring\\u201d \" +\n  \"whereas elsewhere in this manuscript we have \\u201coptimise\\u201d (see for \" +\n  \"example page 26, last two lines) and \\u201cbehaviour\\u201d (see for example, \" +\n  \"page 27, Section 3.2.1, paragraph 1, line 10). Please make sure that \" +\n  \"\\u201cor\\u201d or \\u201cour\\u201d and \\u201cise\\u201d or \\u201cize\\u201d usage is \" +\n  \"standardized throughout for all relevant words. I have no preference for \" +\n  \"which is chosen but consistency is important.\";\n\nconst results1 = body.search(target1, { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  results1.items[0].font.color = \"#C9211E\";\n}\n\n// --- Edit 2: \"... Please check the spelling of 'Murril' ...\" comment -\nconst target2 =\n  \"Please check the spelling of \\u201cMurril\\u201d (or \\u201cMurrill\\u201d) \\u2013 \" +\n  \"you have spelled them differently in various places. On page 155, you \" +\n  \"have two references: \\u201cMurril PW (1967)\\u201d and \\u201cRovira AJA, Murrill \" +\n  \"PW, Smith CL (1969)\\u201d and there are various mentions in the text (a \" +\n  \"method of tuning among them). Presumably this is the same individual \" +\n  \"who only spells his name in one way or the other.\";\n\nconst results2 = body.search(target2, { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  results2.items[0].font.color = \"#C9211E\";\n}\n\nawait context.sync();\n", "ps1": "# Apply the two red-colored annotation edits requested by Oliver.\n# Both edits use Find.Execute to select the \"comment\" text that follows\n# the existing \"In General:\" lead-in run (leaving the separating space\n# in its own, unformatted run, matching the reviewer's manual edit) and\n# then set Font.Color on the exact matched range so Word (re)splits the\n# underlying runs as needed. Any existing highlight on sub-runs (e.g.\n# \"Murril\"/\"Murrill\"/\"PW\") is preserved because Color and Highlight are\n# independent character formatting properties.\n\n$d = $word.ActiveDocument\n\n# Word/VBA Font.Color is a BGR-packed long (R + G*256 + B*65536), not a\n# straight hex RGB value, so build 0xC9211E from its components.\n$redColor = 0xC9 + (0x21 * 256) + (0x1E * 65536)\n\n# --- Edit 1: \"... on page 56 we find ...\" comment -----------------\n$target1 = \"on page 56 we find \u201coptimization\u201d and \u201cneighboring\u201d whereas elsewhere in this manuscript we have \u201coptimise\u201d (see for example page 26, last two lines) and \u201cbehaviour\u201d (see for example, page 27, Section 3.2.1, paragraph 1, line 10). Please make sure that \u201cor\u201d or \u201cour\u201d and \u201cise\u201d or \u201cize\u201d usage is standardized throughout for all relevant words. I have no preference for which is chosen but consistency is important.\"\n\n$rng1 = $d.Content\nif ($rng1.Find.Execute($target1)) {\n    $rng1.Font.Color = $redColor\n}\n\n# --- Edit 2: \"... Please check the spelling of 'Murril' ...\" comment -\n$target2 = \"Please check the spelling of \u201cMurril\u201d (or \u201cMurrill\u201d) \u2013 you have spelled them differently in various places. On page 155, you have two references: \u201cMurril PW (1967)\u201d and \u201cRovira AJA, Murrill PW, Smith CL (1969)\u201d and there are various mentions in the text (a method of tuning among them). Presumably this is the same individual who only spells his name in one way or the other.\"\n\n$rng2 = $d.Content\nif ($rng2.Find.Execute($target2)) {\n    $rng2.Font.Color = $redColor\n}\n"}
